$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RVL")

# Remove the "Map Range" parameter rows (fromRow/fromCol/toRow/toCol),
# which shifts everything below them up by 4 rows.
$ws.Rows("9:12").Delete()

# The "Action Functions" block (EsLaunch / EpChangeCompany / EsOpenModule /
# EpClickRibbon) now sits at rows 11-14. Re-point these rows at the new "EP"
# global object, and rename the two EP actions to their short forms.
$ws.Range("C11").Value = "EP"
$ws.Range("C12").Value = "EP"
$ws.Range("D12").Value = "ChangeCompany"
$ws.Range("C13").Value = "EP"
$ws.Range("C14").Value = "EP"
$ws.Range("D14").Value = "ClickRibbon"
